$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3,4,5,7,8,9,10,11,12 have their Fecha/Calidad/Volumen/Precio columns
# (D, L, M, N, O, P, S) rotated among each other (row 6 is untouched).
# Capture the "before" values first so the rotation can be applied safely
# without one row's new value clobbering data still needed by another row.

$rows = @(3, 4, 5, 7, 8, 9, 10, 11, 12)
$cols = @("D", "L", "M", "N", "O", "P", "S")

$original = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowData
}

# Destination row -> source row (i.e. destination row receives the
# original data that used to live in the source row).
$mapping = @{
    3  = 8
    4  = 11
    5  = 12
    7  = 10
    8  = 4
    9  = 5
    10 = 9
    11 = 7
    12 = 3
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcData = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcData[$c]
    }
}
